$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the statistical test changed from ANOVA to Kruskal-Wallis,
# with updated stat and p-value results (uniformity test / tester fix).
$ws.Range("D2").Value = "Kruskal-Wallis"
$ws.Range("E2").Value = 15.41868384604637
$ws.Range("F2").Value = 0.5653518641114013
$ws.Range("D3").Value = "Kruskal-Wallis"
$ws.Range("E3").Value = 16.97034382861938
$ws.Range("F3").Value = 0.4563770906250505
$ws.Range("D4").Value = "Kruskal-Wallis"
$ws.Range("E4").Value = 13.91313131313132
$ws.Range("F4").Value = 0.6732412363272156
$ws.Range("D5").Value = "Kruskal-Wallis"
$ws.Range("E5").Value = 16.49898989898995
$ws.Range("F5").Value = 0.4887769044145326
$ws.Range("D9").Value = "Kruskal-Wallis"
$ws.Range("E9").Value = 13.83472471060628
$ws.Range("F9").Value = 0.6787646812994838
$ws.Range("D10").Value = "Kruskal-Wallis"
$ws.Range("E10").Value = 13.6437710437711
$ws.Range("F10").Value = 0.6921379188669889
$ws.Range("D11").Value = "Kruskal-Wallis"
$ws.Range("E11").Value = 8.135353535353545
$ws.Range("F11").Value = 0.963566203708729
$ws.Range("D12").Value = "Kruskal-Wallis"
$ws.Range("E12").Value = 12.70639730639729
$ws.Range("F12").Value = 0.7556082428860643
$ws.Range("D16").Value = "Kruskal-Wallis"
$ws.Range("E16").Value = 6.730222866000878
$ws.Range("F16").Value = 0.9867805605494859
$ws.Range("D17").Value = "Kruskal-Wallis"
$ws.Range("E17").Value = 12.96480394399127
$ws.Range("F17").Value = 0.7385389937130951
$ws.Range("D18").Value = "Kruskal-Wallis"
$ws.Range("E18").Value = 3.359595959595993
$ws.Range("F18").Value = 0.9998446593377505
$ws.Range("D19").Value = "Kruskal-Wallis"
$ws.Range("E19").Value = 20.4639730639731
$ws.Range("F19").Value = 0.2511879998136732
$ws.Range("D23").Value = "Kruskal-Wallis"
$ws.Range("E23").Value = 15.30841750841753
$ws.Range("F23").Value = 0.5732772659217072
$ws.Range("D24").Value = "Kruskal-Wallis"
$ws.Range("E24").Value = 13.92121212121216
$ws.Range("F24").Value = 0.6726709897822611
$ws.Range("D25").Value = "Kruskal-Wallis"
$ws.Range("E25").Value = 12.20808080808081
$ws.Range("F25").Value = 0.7873765668019929
$ws.Range("D26").Value = "Kruskal-Wallis"
$ws.Range("E26").Value = 12.15420875420875
$ws.Range("F26").Value = 0.7907106109798886
$ws.Range("D30").Value = "Kruskal-Wallis"
$ws.Range("E30").Value = 6.360269360269371
$ws.Range("F30").Value = 0.9904193734119163
$ws.Range("D31").Value = "Kruskal-Wallis"
$ws.Range("E31").Value = 3.620875420875421
$ws.Range("F31").Value = 0.9997381806747955
$ws.Range("D32").Value = "Kruskal-Wallis"
$ws.Range("E32").Value = 2.050505050505052
$ws.Range("F32").Value = 0.9999958389324287
$ws.Range("D33").Value = "Kruskal-Wallis"
$ws.Range("E33").Value = 2.979797979798008
$ws.Range("F33").Value = 0.9999337916661879
$ws.Range("D37").Value = "Kruskal-Wallis"
$ws.Range("E37").Value = 10.76700336700335
$ws.Range("F37").Value = 0.8684166538117963
$ws.Range("D38").Value = "Kruskal-Wallis"
$ws.Range("E38").Value = 6.141648242738463
$ws.Range("F38").Value = 0.9921822497078893
$ws.Range("D39").Value = "Kruskal-Wallis"
$ws.Range("E39").Value = 2.190572390572413
$ws.Range("F39").Value = 0.9999931415624029
$ws.Range("D40").Value = "Kruskal-Wallis"
$ws.Range("E40").Value = 3.661279461279491
$ws.Range("F40").Value = 0.9997173307080688
